$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matching original inlineStr cells)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6:D15").NumberFormat = "@"
$ws.Range("D17:D23").NumberFormat = "@"
$ws.Range("D25:D30").NumberFormat = "@"
$ws.Range("D32:D51").NumberFormat = "@"

$ws.Range("D2").Value = '22.428.37'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.563.28'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '288.43'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").Value = '0.3669'
$ws.Range("E7").Value = '  -2.00%  '
$ws.Range("D8").Value = '49.84'
$ws.Range("E8").Value = '  +1.10%  '
$ws.Range("D9").Value = '0.3360'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '1.132'
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("D11").Value = '0.07465'
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '20.91'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").Value = '5.965'
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").Value = '6.934'
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("D16").Value = '1.565.16'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '0.00001107'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '90.00'
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Value = '0.06737'
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = '6.344'
$ws.Range("E21").Value = '  +3.13%  '
$ws.Range("D22").Value = '16.13'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '12.01'
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").Value = '22.408.95'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '2.394'
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("D26").Value = '2.613'
$ws.Range("E26").Value = '  +3.76%  '
$ws.Range("D27").Value = '19.72'
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").Value = '149.09'
$ws.Range("E28").Value = '  +1.28%  '
$ws.Range("D29").Value = '5.046'
$ws.Range("E29").Value = '  +0.98%  '
$ws.Range("D30").Value = '123.67'
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("D31").Value = '1.738.48'
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").Value = '1.050'
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("D33").Value = '2.020'
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("D34").Value = '6.117'
$ws.Range("E34").Value = '  +3.95%  '
$ws.Range("D35").Value = '9.604'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").Value = '0.08273'
$ws.Range("E36").Value = '  -1.67%  '
$ws.Range("D37").Value = '0.02428'
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").Value = '1.331'
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("D39").Value = '0.2258'
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("D40").Value = '0.06405'
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").Value = '5.316'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").Value = '11.10'
$ws.Range("E42").Value = '  -1.98%  '
$ws.Range("D43").Value = '0.6133'
$ws.Range("E43").Value = '  -1.04%  '
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").Value = '13.80'
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("D46").Value = '3.762'
$ws.Range("E46").Value = '  -1.21%  '
$ws.Range("D47").Value = '0.5758'
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").Value = '2.028'
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("D49").Value = '125.51'
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").Value = '1.218'
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("D51").Value = '0.07323'
$ws.Range("E51").Value = '  +0.35%  '
